# Fix typo in question name: "certifcat-indigence" -> "certificat-indigence"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("B10").Value = "certificat-indigence"

# Move the active selection, matching the saved state after the edit
$ws.Range("C12").Select()
